$d = $word.ActiveDocument

# --- 1) Body paragraph: merge "set) in the range of 1" + <bookmark> + " to 1000"
#        into a single run "set) in the range of 1 to 1000" (bookmark removed).
#        A two-phase FormattedText replace avoids the engine's run-merge cascade
#        bleeding into the neighbouring "use a " / "0.  The program" runs.
$full = $d.Content.Text
$i1 = $full.IndexOf("set) in the range of 1")
$i2 = $full.IndexOf(" to 1000")
$r = $d.Range($i1, $i2 + 8)
$ft = $r.FormattedText
$ft.Text = "set) in the range of 1 to 0100"
$r.FormattedText = $ft

$full2 = $d.Content.Text
$j1 = $full2.IndexOf("set) in the range of 1 to 0100")
$r2 = $d.Range($j1, $j1 + 30)
$ft2 = $r2.FormattedText
$ft2.Text = "set) in the range of 1 to 1000"
$r2.FormattedText = $ft2

# --- 2) Title paragraph: "Lab Exercise 12/12/2019" -> "Lab Exercise 1/5/2021",
#        with the (now-relocated) _GoBack bookmark landing between "1/5" and "/2021",
#        splitting the run in two just like Word's own last-edit-position bookmark.
$full3 = $d.Content.Text
$i3 = $full3.IndexOf("Lab Exercise 12/12/2019")
$len3 = "Lab Exercise 12/12/2019".Length
$r3 = $d.Range($i3, $i3 + $len3)
$r3.Text = "Lab Exercise 1/5/2021"

$full4 = $d.Content.Text
$j3 = $full4.IndexOf("Lab Exercise 1/5/2021")
$bmPos = $j3 + "Lab Exercise 1/5".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
